$d = $word.ActiveDocument

$replacements = @(
    @("2023-09-13 Wednesday", "2023-09-14 Thursday"),
    @("31×59=", "40×76="),
    @("53×48=", "65×85="),
    @("92×52=", "47×92="),
    @("65×21=", "46×74="),
    @("33×69=", "78×60="),
    @("70×92=", "52×74="),
    @("33×36=", "40×79="),
    @("79×46=", "77×88="),
    @("37×76=", "33×45="),
    @("50×81=", "36×57="),
    @("35×81=", "61×66="),
    @("28×94=", "27×77="),
    @("22×64=", "70×36="),
    @("31×60=", "22×50="),
    @("50×77=", "96×90="),
    @("95×57=", "81×53="),
    @("92×33=", "56×64="),
    @("12×16=", "85×28="),
    @("80×12=", "40×11="),
    @("12×26=", "95×32="),
    @("79×20=", "12×65="),
    @("48×19=", "40×31="),
    @("78×39=", "61×42="),
    @("85×98=", "86×41="),
    @("42×65=", "54×99=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new, 2)
}

$d.Save()
